# Add support for header bidding
# - Adds a new shared string (the header-bidding response JSON) and
#   assigns it to cells F4:F11 (previously empty).
# - Grows rows 4-11 to the maximum row height (409.5) to accommodate the
#   new, much longer response text.
# - Updates the sheet's active selection to F10:F11 (active cell F11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$hbResponse = '[{"ad":"<html><body><script type=\"text\/javascript\">var ut_ju = ''http://stgads.undertone.com/aj'';ut = new Object();ut.bidid=''2ae64e05296983'';ut.bannerid=1197504;ut.zoneid=173879;ut.hbprice=2.08;ut.width=1;ut.height=1;ut.adaptor=''prebid'';<\/script><script type=\"text/javascript\" src=\"http://cdn.undertone.com/js/ajs.js\"><\/script><\/script><\/body><\/html>","publisherId": 3470,"bidRequestId": "2ae64e05296983","placementId": "10433394","adId": 1197504,"campaignId": 297790,"height": 1,"width": 1,"ttl": 700,"currency": "USD","cpm": 2.08,"adaptor": "prebid","netRevenue": "true"}]'

foreach ($r in 4..11) {
    $ws.Cells.Item($r, 6).Value = $hbResponse
    $ws.Rows.Item($r).RowHeight = 409.5
}

$ws.Activate()
$ws.Range("F10:F11").Select()
